# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - rows 3,4,5 column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 303
$wsExhibit.Range("F4").Value = 1307
$wsExhibit.Range("F5").Value = 83

# Sheet "全部类型" (All types) - rows 4,5,6 column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 303
$wsAll.Range("F5").Value = 1307
$wsAll.Range("F6").Value = 83

$wb.Save()
